# Apply updated cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "66.261.91"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "3.566.56"
$ws.Range("E3").Value = "  +4.60%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.89"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.25"
$ws.Range("E6").Value = "  +1.89%  "

$ws.Range("D7").Value = "3.565.36"
$ws.Range("E7").Value = "  +4.60%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("E9").Value = "  +3.40%  "

$ws.Range("E10").Value = "  +1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.03"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").Value = "4.170.56"
$ws.Range("E13").Value = "  +4.60%  "

$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").Value = "3.563.76"
$ws.Range("E16").Value = "  +4.47%  "

$ws.Range("D17").Value = "66.363.96"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.49"
$ws.Range("E19").Value = "  +9.99%  "

$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.97"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.34"
$ws.Range("E22").Value = "  +3.52%  "

$ws.Range("E23").Value = "  +4.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.73"
$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("D25").Value = "3.708.00"
$ws.Range("E25").Value = "  +4.52%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  +7.29%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  +3.82%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.160"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").Value = "3.559.69"
$ws.Range("E34").Value = "  +4.56%  "

$ws.Range("E35").Value = "  +3.72%  "

$ws.Range("E36").Value = "  +3.27%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  +4.08%  "

$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.77"
$ws.Range("E41").Value = "  +1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0856"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  +2.73%  "

$ws.Range("E44").Value = "  +2.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.97"
$ws.Range("E45").Value = "  +1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.06"
$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("E47").Value = "  +3.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.02"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("E51").Value = "  +12.84%  "
